$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(3,6,10,11,12,13,14,15,17,18,19,20,30,33,37,38,39,40,41,42,44,45,46,47,57,60,64,65,66,67,68,69,71,72,73,74,86,87,88,89,93,95,96,97,99,112,113,114,115,119,121,122,123,125,138,139,140,141,145,147,148,149,151)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    if ($cell.Text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
